# ---------------------------------------------------------------------------
# Blog_stack_value_results.xlsx edit
#
# Commit message: "When looking up MLB player ID by display name, ignore
# ' Jr'" -- the workbook itself is the author's weekly "did the stack hit"
# blog-post generator: the "Current" sheet holds this week's two stacks
# (Boston Red Sox hitters, St Louis Cardinals righties), "RG table" turns
# those into pipe-delimited blog text, and "Season Log" is the running
# history of every stack ever posted. Populate this week's numbers and
# append the corresponding Season Log rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsCurrent = $wb.Worksheets.Item("Current")
$wsRG      = $wb.Worksheets.Item("RG table")
$wsSeason  = $wb.Worksheets.Item("Season Log")

# ------------------------------------------------------------------
# 1. "Current" sheet -- fill in this week's two stacks
# ------------------------------------------------------------------

# Results of the two stacks (row 8, under the value column). These are
# written first -- while B3:C6/G3:H6 are still blank and D7/I7 are still
# #DIV/0! -- so that the "RG table" CONCATENATE formulas (which stitch
# together Current!D8 / Current!I8 along with the rest of the stack row)
# pick them up cleanly as the rest of the numbers are filled in below,
# rather than leaving a stale cached blank for D8/I8 in those formulas.
$wsCurrent.Range("D8").Value2 = "Failure"
$wsCurrent.Range("I8").Value2 = "Success"

# Header titles for the two stacks in the top block (row 1)
$wsCurrent.Range("A1").Value2 = "Boston Red Sox hitters (FD, DK)"
$wsCurrent.Range("F1").Value2 = "St Louis Cardinals righties (FD, DK)"

# Stack 1 -- Boston Red Sox hitters (columns A:D)
$wsCurrent.Range("A3").Value2 = "Betts"
$wsCurrent.Range("B3").Value2 = 4300
$wsCurrent.Range("C3").Value2 = 25.1

$wsCurrent.Range("A4").Value2 = "Martinez"
$wsCurrent.Range("B4").Value2 = 4400
$wsCurrent.Range("C4").Value2 = 3

$wsCurrent.Range("A5").Value2 = "Vazquez"
$wsCurrent.Range("B5").Value2 = 3000
$wsCurrent.Range("C5").Value2 = 12.2

# Row 6 for this stack is left blank, same as before the edit.

# Stack 2 -- St Louis Cardinals righties (columns F:I)
$wsCurrent.Range("F3").Value2 = "Edman"
$wsCurrent.Range("G3").Value2 = 2800
$wsCurrent.Range("H3").Value2 = 15.5

$wsCurrent.Range("F4").Value2 = "Martinez"
$wsCurrent.Range("G4").Value2 = 2700
$wsCurrent.Range("H4").Value2 = 18.7

$wsCurrent.Range("F5").Value2 = "Goldschmidt"
$wsCurrent.Range("G5").Value2 = 3100
$wsCurrent.Range("H5").Value2 = 28.2

$wsCurrent.Range("F6").Value2 = "O'Neill"
$wsCurrent.Range("G6").Value2 = 2300
$wsCurrent.Range("H6").Value2 = 12.2

# ------------------------------------------------------------------
# 2. "Season Log" sheet -- backfill one missing result + append the
#    new history rows for this week's (and the past week's) stacks.
# ------------------------------------------------------------------

# Row 198 was entered with only a date/site; fill in the rest of it.
$wsSeason.Range("C198").Value2 = "Milwaukee Brewers hitters (FD, DK)"
$wsSeason.Range("D198").Value2 = 4.37
$wsSeason.Range("E198").Value2 = "Success"

# New rows appended to the log.
$seasonRows = @(
    @(43648, "RG",        "Tampa Bay Rays lefties (FD, DK)",      3.02, "Failure"),
    @(43648, "RG",        "Colorado Rockies hitters (FD, DK)",    3.22, "Failure"),
    @(43648, "RG",        "Los Angeles Dodgers hitters (FD, DK)", 2.09, "Failure"),
    @(43649, "RG",        "Cincinnati Reds lefties (FD, DK)",     1.27, "Failure"),
    @(43649, "RG",        "Colorado Rockies righties (FD, DK)",   0.84, "Failure"),
    @(43649, "RG",        "Los Angeles Angels hitters (FD, DK)",  5.16, "Success"),
    @(43650, "Draftshot", "Los Angeles Dodgers lefties (FD, DK)", 3.6,  "Failure"),
    @(43651, "RG",        "Boston Red Sox hitters (FD, DK)",      3.44, "Failure"),
    @(43651, "RG",        "St Louis Cardinals righties (FD, DK)", 6.84, "Success")
)

$r = 201
foreach ($row in $seasonRows) {
    $wsSeason.Cells.Item($r, 1).Value2 = $row[0]
    $wsSeason.Cells.Item($r, 2).Value2 = $row[1]
    $wsSeason.Cells.Item($r, 3).Value2 = $row[2]
    $wsSeason.Cells.Item($r, 4).Value2 = $row[3]
    $wsSeason.Cells.Item($r, 5).Value2 = $row[4]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 3. View / selection state to match the author's saved workbook:
#    - "Current": selection moved off F18 onto the merged F1:I1 title
#      cell, and it's no longer the active sheet.
#    - "RG table": A1:A2 selected (the two freshly-generated blog
#      blurbs).
#    - "Season Log": becomes the active sheet/tab, with the selection
#      sitting on the newly appended last row.
# ------------------------------------------------------------------

$wsCurrent.Range("F1:I1").Select()
$wsRG.Range("A1:A2").Select()

$wsSeason.Activate()
$wsSeason.Range("F209").Select()
